# Applies the "meter_per_cell / cell_per_axis / min_obstacle_height" relabeling
# edits to the single-slide figure deck.
#
# Shape map (Slide 1, Shapes.Item index -> cNvPr id / name / original text):
#   7  -> id 47 "TextShape 7"  "m_per_cell"
#   10 -> id 50 "TextShape 10" "grid_dim"
#   15 -> id 55 "TextShape 15" "h > height_diff_threshold"
#   23 -> id 63 "TextShape 22" "h < height_diff_threshold"
#   24 -> id 25 "TextShape 15" "h > height_diff_threshold"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) Shape 7: "m_per_cell" -> "meter_per_cell", plus a resize/reposition of
#    the textbox (it grows to the left and gets wider).
# ---------------------------------------------------------------------------
$sh7 = $s.Shapes.Item(7)
$sh7.TextFrame.TextRange.Text = "meter_per_cell"
$sh7.Left = 424.148071
$sh7.Width = 63.411063

# ---------------------------------------------------------------------------
# 2) Shape 10: "grid_dim" -> "cell" + "_per_axis" (two runs).
# ---------------------------------------------------------------------------
$sh10 = $s.Shapes.Item(10)
$tr10 = $sh10.TextFrame.TextRange
$tr10.Text = "cell_per_axis"
# Re-touch the second run's formatting so the text engine splits it away from
# the first one (both keep the shape's existing Arial/6pt/-1 spacing style).
$tr10.Characters(5, 9).Font.Bold = 0

# ---------------------------------------------------------------------------
# 3) Shape 15: "h > height_diff_threshold" -> "h >" / " " / "min_obstacle_height"
#    (three runs).
# ---------------------------------------------------------------------------
$sh15 = $s.Shapes.Item(15)
$tr15 = $sh15.TextFrame.TextRange
$tr15.Characters(5, 21).Text = "min_obstacle_height"
$tr15.Characters(4, 1).Font.Name = "Arial"

# ---------------------------------------------------------------------------
# 4) Shape 23 (TextShape 22): "h < height_diff_threshold" ->
#    "h < " / "min_obstacle_height" (two runs).
# ---------------------------------------------------------------------------
$sh23 = $s.Shapes.Item(23)
$tr23 = $sh23.TextFrame.TextRange
$tr23.Characters(5, 21).Text = "min_obstacle_height"

# ---------------------------------------------------------------------------
# 5) Shape 24 (second "TextShape 15"): same change as shape 15.
# ---------------------------------------------------------------------------
$sh24 = $s.Shapes.Item(24)
$tr24 = $sh24.TextFrame.TextRange
$tr24.Characters(5, 21).Text = "min_obstacle_height"
$tr24.Characters(4, 1).Font.Name = "Arial"
